# #5: cash & deposit done
# Sheet "存款" (deposits): turn row 1 into a proper header row and add the
# per-row provenance columns (property_category..index) that already exist
# on the other sheets (土地/建物/汽車).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Row 1: header labels -------------------------------------------------
# B1/C1/D1 used to just repeat row 2's data; they become column headers.
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Style the new header cells like the rest of row 1 (bold, centered, boxed).
$hdr = $ws.Range("G1:M1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# The "date" column holds a literal text value ("2011-12-26"), not a real
# date -- format as text first so Excel doesn't coerce it into a serial.
$ws.Range("I2:I8").NumberFormat = "@"

# --- Rows 2-8: fill in the new provenance columns -------------------------
$indexes = @(48, 49, 50, 51, 52, 53, 54)

for ($i = 0; $i -lt $indexes.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 7).Value  = "deposit"        # G: property_category
    $ws.Cells.Item($r, 8).Value  = "normal"          # H: category
    $ws.Cells.Item($r, 9).Value  = "2011-12-26"      # I: date
    $ws.Cells.Item($r, 10).Value = "潘孟安"           # J: legislator_name
    $ws.Cells.Item($r, 11).Value = 1376              # K: legislator_id
    $ws.Cells.Item($r, 12).Value = "tmp6a821"        # L: source_file
    $ws.Cells.Item($r, 13).Value = $indexes[$i]      # M: index
}

Write-Output "sheet4 (deposits) header + provenance columns updated"
